$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 73819.766
$ws.Range("J3").Value = 73819.766
$ws.Range("L3").Value = 73819.766
$ws.Range("N3").Value = -74047.766
$ws.Range("H5").Value = 1333412.9
$ws.Range("I5").Value = 1818269
$ws.Range("J5").Value = 58.75
$ws.Range("K5").Value = 1818269
$ws.Range("L5").Value = 58.75
$ws.Range("M5").Value = -1818154
$ws.Range("N5").Value = -288.75
$ws.Range("H33").Value = 562.2
$ws.Range("I33").Value = 251.70589
$ws.Range("K33").Value = 251.70589
$ws.Range("M33").Value = -22.70589000000001
$ws.Range("H75").Value = 62000
$ws.Range("J75").Value = 62000
$ws.Range("L75").Value = 62000
$ws.Range("N75").Value = -63872
$ws.Range("H78").Value = 62000
$ws.Range("J78").Value = 62000
$ws.Range("L78").Value = 186000
$ws.Range("N78").Value = -195360
$ws.Range("H102").Value = 73819.766
$ws.Range("J102").Value = 73819.766
$ws.Range("L102").Value = 73819.766
$ws.Range("N102").Value = -80309.766
$ws.Range("H112").Value = 1688.0834
$ws.Range("J112").Value = 1669.5454
$ws.Range("L112").Value = 5008.6362
$ws.Range("N112").Value = -7224.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 801.9841300000001
$ws.Range("I2").Value = 744.62067
$ws.Range("K2").Value = 744.62067
$ws.Range("M2").Value = -631.62067
$ws.Range("H32").Value = 3855.0425
$ws.Range("I32").Value = 1989.6097
$ws.Range("K32").Value = 1989.6097
$ws.Range("M32").Value = -1702.6097
$ws.Range("H63").Value = 2052.3872
$ws.Range("I63").Value = 1658.5
$ws.Range("K63").Value = 1658.5
$ws.Range("M63").Value = -972.5
$ws.Range("H66").Value = 2052.3872
$ws.Range("I66").Value = 1658.5
$ws.Range("K66").Value = 8292.5
$ws.Range("M66").Value = -4860.5
$ws.Range("H116").Value = 801.9841300000001
$ws.Range("I116").Value = 744.62067
$ws.Range("K116").Value = 744.62067
$ws.Range("M116").Value = 1549.37933
$ws.Range("H122").Value = 2641.8
$ws.Range("I122").Value = 3005.5
$ws.Range("J122").Value = 2399.3333
$ws.Range("K122").Value = 9016.5
$ws.Range("L122").Value = 7197.999899999999
$ws.Range("M122").Value = -6566.5
$ws.Range("N122").Value = -12097.9999
$ws.Range("H132").Value = 2374.7896
$ws.Range("I132").Value = 1730.5
$ws.Range("J132").Value = 2450.5881
$ws.Range("K132").Value = 5191.5
$ws.Range("L132").Value = 7351.7643
$ws.Range("M132").Value = -2661.5
$ws.Range("N132").Value = -12411.7643

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 801.9841300000001
$ws.Range("I3").Value = 744.62067
$ws.Range("K3").Value = 744.62067
$ws.Range("M3").Value = -630.62067
$ws.Range("H29").Value = 4666.3335
$ws.Range("I29").Value = 4666.3335
$ws.Range("K29").Value = 4666.3335
$ws.Range("M29").Value = -4377.3335
$ws.Range("H64").Value = 578.8
$ws.Range("I64").Value = 765
$ws.Range("K64").Value = 765
$ws.Range("M64").Value = -540
$ws.Range("H67").Value = 578.8
$ws.Range("I67").Value = 765
$ws.Range("K67").Value = 765
$ws.Range("M67").Value = 15
$ws.Range("H82").Value = 18582.223
$ws.Range("I82").Value = 12464.286
$ws.Range("J82").Value = 39995
$ws.Range("K82").Value = 12464.286
$ws.Range("L82").Value = 39995
$ws.Range("M82").Value = -12081.286
$ws.Range("N82").Value = -40761
$ws.Range("H85").Value = 18582.223
$ws.Range("I85").Value = 12464.286
$ws.Range("J85").Value = 39995
$ws.Range("K85").Value = 12464.286
$ws.Range("L85").Value = 39995
$ws.Range("M85").Value = -11138.286
$ws.Range("N85").Value = -42647
$ws.Range("H97").Value = 40117.75
$ws.Range("I97").Value = 33333.332
$ws.Range("K97").Value = 33333.332
$ws.Range("M97").Value = -32342.332
$ws.Range("H105").Value = 4896.25
$ws.Range("I105").Value = 4896.25
$ws.Range("K105").Value = 4896.25
$ws.Range("M105").Value = -3149.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6092.8667
$ws.Range("I58").Value = 6384.077
$ws.Range("K58").Value = 6384.077
$ws.Range("M58").Value = -6181.077
$ws.Range("H86").Value = 5498.727
$ws.Range("I86").Value = 5397
$ws.Range("K86").Value = 5397
$ws.Range("M86").Value = -4274
$ws.Range("H89").Value = 5498.727
$ws.Range("I89").Value = 5397
$ws.Range("K89").Value = 26985
$ws.Range("M89").Value = -21369
$ws.Range("H103").Value = 6598.4
$ws.Range("I103").Value = 6598.4
$ws.Range("K103").Value = 6598.4
$ws.Range("M103").Value = -5426.4
$ws.Range("H122").Value = 2494.4
$ws.Range("I122").Value = 2887.8572
$ws.Range("J122").Value = 1576.3334
$ws.Range("K122").Value = 8663.571599999999
$ws.Range("L122").Value = 4729.0002
$ws.Range("M122").Value = -6213.571599999999
$ws.Range("N122").Value = -9629.0002
$ws.Range("H136").Value = 6092.8667
$ws.Range("I136").Value = 6384.077
$ws.Range("K136").Value = 19152.231
$ws.Range("M136").Value = -16602.231
$ws.Range("H141").Value = 49997.8
$ws.Range("J141").Value = 56247.75
$ws.Range("L141").Value = 56247.75
$ws.Range("N141").Value = -66607.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 53.375
$ws.Range("J2").Value = 26.666666
$ws.Range("L2").Value = 159.999996
$ws.Range("N2").Value = -385.999996
$ws.Range("H113").Value = 715
$ws.Range("J113").Value = 769.125
$ws.Range("L113").Value = 2307.375
$ws.Range("N113").Value = -6647.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13482.583
$ws.Range("J80").Value = 19398.715
$ws.Range("L80").Value = 19398.715
$ws.Range("N80").Value = -21394.715
$ws.Range("H83").Value = 13482.583
$ws.Range("J83").Value = 19398.715
$ws.Range("L83").Value = 96993.575
$ws.Range("N83").Value = -106977.575
$ws.Range("H93").Value = 23512.777
$ws.Range("I93").Value = 20227.666
$ws.Range("J93").Value = 30083
$ws.Range("K93").Value = 20227.666
$ws.Range("L93").Value = 30083
$ws.Range("M93").Value = -18355.666
$ws.Range("N93").Value = -33827
$ws.Range("H104").Value = 671
$ws.Range("J104").Value = 671
$ws.Range("L104").Value = 671
$ws.Range("N104").Value = -7659
$ws.Range("H122").Value = 2679.4348
$ws.Range("I122").Value = 2660
$ws.Range("J122").Value = 2691.9285
$ws.Range("K122").Value = 7980
$ws.Range("L122").Value = 8075.7855
$ws.Range("M122").Value = -5530
$ws.Range("N122").Value = -12975.7855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2626.4285
$ws.Range("I22").Value = 3247
$ws.Range("J22").Value = 1799
$ws.Range("K22").Value = 3247
$ws.Range("L22").Value = 1799
$ws.Range("M22").Value = -2952
$ws.Range("N22").Value = -2389
$ws.Range("H27").Value = 2626.4285
$ws.Range("I27").Value = 3247
$ws.Range("J27").Value = 1799
$ws.Range("K27").Value = 3247
$ws.Range("L27").Value = 1799
$ws.Range("M27").Value = -3140
$ws.Range("N27").Value = -2013
$ws.Range("H55").Value = 238.52632
$ws.Range("I55").Value = 258.91666
$ws.Range("J55").Value = 203.57143
$ws.Range("K55").Value = 258.91666
$ws.Range("L55").Value = 203.57143
$ws.Range("M55").Value = -85.91665999999998
$ws.Range("N55").Value = -549.57143

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 18000
$ws.Range("J104").Value = 18000
$ws.Range("L104").Value = 18000
$ws.Range("N104").Value = -24988
$ws.Range("H132").Value = 5061.32
$ws.Range("I132").Value = 4993.2085
$ws.Range("K132").Value = 14979.6255
$ws.Range("M132").Value = -12449.6255
$ws.Range("H136").Value = 4432.3335
$ws.Range("I136").Value = 3832.6667
$ws.Range("J136").Value = 5631.6665
$ws.Range("K136").Value = 11498.0001
$ws.Range("L136").Value = 16894.9995
$ws.Range("M136").Value = -8948.000100000001
$ws.Range("N136").Value = -21994.9995
